$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ToolObject")

# Insert a new column before the old "OutputFood" column (H), shifting it to I,
# to make room for the new "Animation" column.
$ws.Columns.Item(8).Insert()

# Header row
$ws.Cells.Item(1, 8).Value = "Animation"
# Type row
$ws.Cells.Item(2, 8).Value = "string"

# Data rows: default Animation value is "Hold", except the cutting board ("Board"/"도마"),
# which uses "Chop".
for ($r = 3; $r -le 17; $r++) {
    $name = $ws.Cells.Item($r, 2).Value()
    if ($name -eq "Board") {
        $ws.Cells.Item($r, 8).Value = "Chop"
    } else {
        $ws.Cells.Item($r, 8).Value = "Hold"
    }
}
